# This script reproduces the "most recent excel dumps" commit:
# a new data row for equipment "TEMU3584963" is inserted as row 9 of
# Sheet1 (pushing the previous rows 9-12 down to rows 10-13).
#
# New row 9 values:
#   A9 = TEMU3584963      (Equipment Number)
#   B9 = MSC CHANNE       (Vessel)
#   C9 = 00904             (Voyage)
#   D9 = DJSEAA3786849    (WONumber)
#   E9 = 7031961991        (ReferenceNumber)
#   F9 = MEDUMM507023     (BOLNumber)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 9; existing rows 9-12 shift down to 10-13.
$ws.Rows.Item(9).Insert()

# A9: new unique value, not purely numeric so it is safely stored as text.
$ws.Range("A9").Value = "TEMU3584963"

# B9: value already exists elsewhere in the sheet ("MSC CHANNE" in B2) -
# copy/paste the value so the existing shared-string text entry is reused
# verbatim instead of typing it again.
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

# C9: value already exists elsewhere in the sheet ("00904" in C2) - reuse it
# the same way. This also preserves the leading zero as text.
$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

# D9: new unique value, not purely numeric so it is safely stored as text.
$ws.Range("D9").Value = "DJSEAA3786849"

# E9: new unique value that is purely numeric digits. Assigning it directly
# would make Excel store it as a number (losing the text semantics used by
# the rest of this column), so instead stage it as text (leading apostrophe
# forces text entry) in an unused scratch cell, copy just the value over to
# E9, then clean the scratch cell back up.
$scratch = $ws.Range("Z1")
$scratch.Value = "'7031961991"
$scratch.Copy()
$ws.Range("E9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Clear()

# F9: new unique value, not purely numeric so it is safely stored as text.
$ws.Range("F9").Value = "MEDUMM507023"
